$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in column A values for rows 2-4 (codAdmin sequence)
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3

# Update the active selection to A4
$ws.Range("A4").Select()

# Update the window height of the workbook view
$excel.ActiveWindow.Height = 9287
